{"js": "// The document body contains a single 20-row x 5-column table of\n// arithmetic expressions such as \"29-8=\". Every cell's expression is\n// being replaced by a new one (same cell position), so we must target\n// cells by their row/column position rather than by searching for the\n// old text (several old expressions, e.g. \"60-9=\", repeat with\n// different replacements elsewhere in the table).\nconst newValues = [\n  \"43+17=\", \"74+6=\", \"12-5=\", \"90-80=\", \"77+2=\",\n  \"7+7=\", \"68-11=\", \"31-3=\", \"45-29=\", \"79-34=\",\n  \"3+56=\", \"90+3=\", \"5+38=\", \"6+70=\", \"60-40=\",\n  \"60-57=\", \"29+56=\", \"80-33=\", \"83-34=\", \"75+17=\",\n  \"23-11=\", \"70-0=\", \"88-3=\", \"76-0=\", \"11+65=\",\n  \"98-25=\", \"60-50=\", \"12-3=\", \"88-45=\", \"49-46=\",\n  \"96+0=\", \"21+26=\", \"91-87=\", \"47+9=\", \"29+24=\",\n  \"23+29=\", \"67+21=\", \"49+29=\", \"31-8=\", \"57+27=\",\n  \"19-16=\", \"24-2=\", \"36+21=\", \"10+35=\", \"6-0=\",\n  \"0+4=\", \"45+39=\", \"72-40=\", \"56-54=\", \"19+57=\",\n  \"38+2=\", \"73+25=\", \"98-65=\", \"24+51=\", \"75-49=\",\n  \"85-22=\", \"37+45=\", \"36+8=\", \"91-59=\", \"12+75=\",\n  \"9+73=\", \"76-71=\", \"71-25=\", \"50-17=\", \"93-77=\",\n  \"74+16=\", \"12+4=\", \"83-5=\", \"92-84=\", \"18+77=\",\n  \"17+56=\", \"50+42=\", \"40-18=\", \"44+17=\", \"85-14=\",\n  \"75+19=\", \"34+56=\", \"6+30=\", \"65+33=\", \"49-18=\",\n  \"43-17=\", \"49+15=\", \"61+3=\", \"82-38=\", \"96-65=\",\n  \"17+59=\", \"64+17=\", \"73-45=\", \"97-0=\", \"83-30=\",\n  \"1+43=\", \"75-5=\", \"5+29=\", \"93-57=\", \"70-14=\",\n  \"44+16=\", \"54-33=\", \"9+25=\", \"73-49=\", \"58-35=\",\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst COLS = 5;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const idx = r * COLS + c;\n    if (idx >= newValues.length) continue;\n    table.getCell(r, c).value = newValues[idx];\n  }\n}\nawait context.sync();\n", "ps1": "# The document body contains a single 20-row x 5-column table of\n# arithmetic expressions such as \"29-8=\". Every cell's expression is\n# replaced by a new one at the same row/column position; several old\n# expressions repeat (e.g. \"60-9=\" appears twice with two different\n# replacements), so cells must be addressed by position, not by\n# searching for the old text.\n$newValues = @(\n    \"43+17=\", \"74+6=\", \"12-5=\", \"90-80=\", \"77+2=\", \"7+7=\", \"68-11=\", \"31-3=\", \"45-29=\", \"79-34=\",\n    \"3+56=\", \"90+3=\", \"5+38=\", \"6+70=\", \"60-40=\", \"60-57=\", \"29+56=\", \"80-33=\", \"83-34=\", \"75+17=\",\n    \"23-11=\", \"70-0=\", \"88-3=\", \"76-0=\", \"11+65=\", \"98-25=\", \"60-50=\", \"12-3=\", \"88-45=\", \"49-46=\",\n    \"96+0=\", \"21+26=\", \"91-87=\", \"47+9=\", \"29+24=\", \"23+29=\", \"67+21=\", \"49+29=\", \"31-8=\", \"57+27=\",\n    \"19-16=\", \"24-2=\", \"36+21=\", \"10+35=\", \"6-0=\", \"0+4=\", \"45+39=\", \"72-40=\", \"56-54=\", \"19+57=\",\n    \"38+2=\", \"73+25=\", \"98-65=\", \"24+51=\", \"75-49=\", \"85-22=\", \"37+45=\", \"36+8=\", \"91-59=\", \"12+75=\",\n    \"9+73=\", \"76-71=\", \"71-25=\", \"50-17=\", \"93-77=\", \"74+16=\", \"12+4=\", \"83-5=\", \"92-84=\", \"18+77=\",\n    \"17+56=\", \"50+42=\", \"40-18=\", \"44+17=\", \"85-14=\", \"75+19=\", \"34+56=\", \"6+30=\", \"65+33=\", \"49-18=\",\n    \"43-17=\", \"49+15=\", \"61+3=\", \"82-38=\", \"96-65=\", \"17+59=\", \"64+17=\", \"73-45=\", \"97-0=\", \"83-30=\",\n    \"1+43=\", \"75-5=\", \"5+29=\", \"93-57=\", \"70-14=\", \"44+16=\", \"54-33=\", \"9+25=\", \"73-49=\", \"58-35=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n$rows = $t.Rows.Count\n\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $idx = (($r - 1) * $cols) + ($c - 1)\n        if ($idx -lt $newValues.Length) {\n            $t.Cell($r, $c).Range.Text = $newValues[$idx]\n        }\n    }\n}\n"}
